# Update the rental report table:
#  - Row 2 becomes what used to be row 4 (Gilberto Gomez / Monteria / Finalizado / Pagado / 370000)
#  - Row 3 keeps its vehicle but shifts to 10/04-15/04, Florencia/Diana, Sin Iniciar, now Pagado, 120000
#  - Row 4 keeps its vehicle but shifts to 20/03-21/03, Monteria/Diana, now Finalizado, 90000
# Dates are stored as plain text in this report, so force text formatting on the
# date cells before assigning them - otherwise Excel would auto-convert the
# "dd/mm/yyyy" strings into date serial numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-TextValue $ws.Range("B2") "10/05/2024"
Set-TextValue $ws.Range("C2") "15/05/2024"
$ws.Range("E2").Value = "Montería, Córdoba"
$ws.Range("F2").Value = "Gilberto Gómez"
$ws.Range("G2").Value = "Finalizado"
$ws.Range("H2").Value = "Pagado"
$ws.Range("I2").Value = 370000
$ws.Range("J2").Value = 55500
$ws.Range("K2").Value = 314500
$ws.Range("L2").Value = 4

# Row 3
Set-TextValue $ws.Range("B3") "10/04/2024"
Set-TextValue $ws.Range("C3") "15/04/2024"
$ws.Range("E3").Value = "Florencia, Caquetá"
$ws.Range("H3").Value = "Pagado"
$ws.Range("I3").Value = 120000
$ws.Range("J3").Value = 18000
$ws.Range("K3").Value = 102000
$ws.Range("L3").Value = 4

# Row 4
Set-TextValue $ws.Range("B4") "20/03/2024"
Set-TextValue $ws.Range("C4") "21/03/2024"
$ws.Range("F4").Value = "Diana Caicedo"
$ws.Range("I4").Value = 90000
$ws.Range("J4").Value = 13500
$ws.Range("K4").Value = 76500
